$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 25,9

$data[0,0] = "model_2_1_24"
$data[0,1] = -0.001236788884735551
$data[0,2] = -0.7474321138812963
$data[0,3] = -0.3349434179765798
$data[0,4] = -0.4304380836533312
$data[0,5] = 1.108074069023132
$data[0,6] = 1.724909663200378
$data[0,7] = 0.2921002507209778
$data[0,8] = 1.050646305084229

$data[1,0] = "model_2_1_22"
$data[1,1] = -0.001236788884735551
$data[1,2] = -0.7474321138812963
$data[1,3] = -0.3349434179765798
$data[1,4] = -0.4304380836533312
$data[1,5] = 1.108074069023132
$data[1,6] = 1.724909663200378
$data[1,7] = 0.2921002507209778
$data[1,8] = 1.050646305084229

$data[2,0] = "model_2_1_21"
$data[2,1] = -0.001236788884735551
$data[2,2] = -0.7474321138812963
$data[2,3] = -0.3349434179765798
$data[2,4] = -0.4304380836533312
$data[2,5] = 1.108074069023132
$data[2,6] = 1.724909663200378
$data[2,7] = 0.2921002507209778
$data[2,8] = 1.050646305084229

$data[3,0] = "model_2_1_20"
$data[3,1] = -0.001236788884735551
$data[3,2] = -0.7474321138812963
$data[3,3] = -0.3349434179765798
$data[3,4] = -0.4304380836533312
$data[3,5] = 1.108074069023132
$data[3,6] = 1.724909663200378
$data[3,7] = 0.2921002507209778
$data[3,8] = 1.050646305084229

$data[4,0] = "model_2_1_19"
$data[4,1] = -0.001236788884735551
$data[4,2] = -0.7474321138812963
$data[4,3] = -0.3349434179765798
$data[4,4] = -0.4304380836533312
$data[4,5] = 1.108074069023132
$data[4,6] = 1.724909663200378
$data[4,7] = 0.2921002507209778
$data[4,8] = 1.050646305084229

$data[5,0] = "model_2_1_18"
$data[5,1] = -0.001236788884735551
$data[5,2] = -0.7474321138812963
$data[5,3] = -0.3349434179765798
$data[5,4] = -0.4304380836533312
$data[5,5] = 1.108074069023132
$data[5,6] = 1.724909663200378
$data[5,7] = 0.2921002507209778
$data[5,8] = 1.050646305084229

$data[6,0] = "model_2_1_17"
$data[6,1] = -0.001236788884735551
$data[6,2] = -0.7474321138812963
$data[6,3] = -0.3349434179765798
$data[6,4] = -0.4304380836533312
$data[6,5] = 1.108074069023132
$data[6,6] = 1.724909663200378
$data[6,7] = 0.2921002507209778
$data[6,8] = 1.050646305084229

$data[7,0] = "model_2_1_16"
$data[7,1] = -0.001236788884735551
$data[7,2] = -0.7474321138812963
$data[7,3] = -0.3349434179765798
$data[7,4] = -0.4304380836533312
$data[7,5] = 1.108074069023132
$data[7,6] = 1.724909663200378
$data[7,7] = 0.2921002507209778
$data[7,8] = 1.050646305084229

$data[8,0] = "model_2_1_23"
$data[8,1] = -0.001236788884735551
$data[8,2] = -0.7474321138812963
$data[8,3] = -0.3349434179765798
$data[8,4] = -0.4304380836533312
$data[8,5] = 1.108074069023132
$data[8,6] = 1.724909663200378
$data[8,7] = 0.2921002507209778
$data[8,8] = 1.050646305084229

$data[9,0] = "model_2_1_15"
$data[9,1] = 0.00050779844049198
$data[9,2] = -0.7350279223295109
$data[9,3] = -0.3360459470563311
$data[9,4] = -0.4217685356563776
$data[9,5] = 1.106143355369568
$data[9,6] = 1.712665200233459
$data[9,7] = 0.2923415303230286
$data[9,8] = 1.044278621673584

$data[10,0] = "model_2_1_14"
$data[10,1] = 0.005154245936211277
$data[10,2] = -0.7023587301430967
$data[10,3] = -0.3372641131688383
$data[10,4] = -0.398693882342235
$data[10,5] = 1.101001143455505
$data[10,6] = 1.68041729927063
$data[10,7] = 0.2926080524921417
$data[10,8] = 1.02733051776886

$data[11,0] = "model_2_1_13"
$data[11,1] = 0.01200897260509437
$data[11,2] = -0.6555721394099363
$data[11,3] = -0.3316818994708288
$data[11,4] = -0.3646215163961446
$data[11,5] = 1.093415021896362
$data[11,6] = 1.634233593940735
$data[11,7] = 0.291386604309082
$data[11,8] = 1.002304553985596

$data[12,0] = "model_2_1_12"
$data[12,1] = 0.02185473053607989
$data[12,2] = -0.583969304067768
$data[12,3] = -0.3451896564012702
$data[12,4] = -0.3155718632502991
$data[12,5] = 1.082518577575684
$data[12,6] = 1.56355357170105
$data[12,7] = 0.2943422496318817
$data[12,8] = 0.966278076171875

$data[13,0] = "model_2_1_11"
$data[13,1] = 0.03205623498687216
$data[13,2] = -0.5101413073049703
$data[13,3] = -0.3567383664693902
$data[13,4] = -0.2646623515686446
$data[13,5] = 1.071228504180908
$data[13,6] = 1.490677237510681
$data[13,7] = 0.2968692481517792
$data[13,8] = 0.928885281085968

$data[14,0] = "model_2_1_10"
$data[14,1] = 0.04026240634942013
$data[14,2] = -0.4474204908867878
$data[14,3] = -0.3773006108488905
$data[14,4] = -0.2229190227406044
$data[14,5] = 1.062146663665771
$data[14,6] = 1.428764820098877
$data[14,7] = 0.3013684749603271
$data[14,8] = 0.8982251286506653

$data[15,0] = "model_2_1_9"
$data[15,1] = 0.048029950845634
$data[15,2] = -0.3920016527940344
$data[15,3] = -0.3783740414318884
$data[15,4] = -0.1836393131046228
$data[15,5] = 1.053550362586975
$data[15,6] = 1.374060273170471
$data[15,7] = 0.301603376865387
$data[15,8] = 0.8693745136260986

$data[16,0] = "model_2_1_8"
$data[16,1] = 0.05832945701982994
$data[16,2] = -0.321727537121409
$data[16,3] = -0.3660651061498115
$data[16,4] = -0.1319135468521053
$data[16,5] = 1.042151808738708
$data[16,6] = 1.304691910743713
$data[16,7] = 0.2989100515842438
$data[16,8] = 0.8313822746276855

$data[17,0] = "model_2_1_7"
$data[17,1] = 0.07040874821892629
$data[17,2] = -0.2364966727899993
$data[17,3] = -0.3666171432725018
$data[17,4] = -0.07134897579323285
$data[17,5] = 1.028783559799194
$data[17,6] = 1.220559597015381
$data[17,7] = 0.2990308403968811
$data[17,8] = 0.7868980765342712

$data[18,0] = "model_2_1_6"
$data[18,1] = 0.106608855335275
$data[18,2] = 0.01530364515557137
$data[18,3] = -0.361239307556261
$data[18,4] = 0.1085584455720723
$data[18,5] = 0.9887207746505737
$data[18,6] = 0.9720046520233154
$data[18,7] = 0.2978541254997253
$data[18,8] = 0.6547573208808899

$data[19,0] = "model_2_1_4"
$data[19,1] = 0.1069855247848572
$data[19,2] = 0.01990723719006282
$data[19,3] = -0.3667089802721606
$data[19,4] = 0.1110615562125522
$data[19,5] = 0.9883038401603699
$data[19,6] = 0.9674603939056396
$data[19,7] = 0.2990509271621704
$data[19,8] = 0.652918815612793

$data[20,0] = "model_2_1_3"
$data[20,1] = 0.1071647629531325
$data[20,2] = 0.02579804654030382
$data[20,3] = -0.3786472967643741
$data[20,4] = 0.1135814389178836
$data[20,5] = 0.9881055951118469
$data[20,6] = 0.9616455435752869
$data[20,7] = 0.3016631603240967
$data[20,8] = 0.6510679721832275

$data[21,0] = "model_2_1_2"
$data[21,1] = 0.1081881867968117
$data[21,2] = 0.04865009374907447
$data[21,3] = -0.4456230786323614
$data[21,4] = 0.1204543361841132
$data[21,5] = 0.9869728684425354
$data[21,6] = 0.9390881061553955
$data[21,7] = 0.3163181841373444
$data[21,8] = 0.6460198760032654

$data[22,0] = "model_2_1_5"
$data[22,1] = 0.1094817261893443
$data[22,2] = 0.03454083708275757
$data[22,3] = -0.3550196428073349
$data[22,4] = 0.1231168375861197
$data[22,5] = 0.9855413436889648
$data[22,6] = 0.9530154466629028
$data[22,7] = 0.2964931726455688
$data[22,8] = 0.6440643668174744

$data[23,0] = "model_2_1_1"
$data[23,1] = 0.110130878546485
$data[23,2] = 0.08345253357714544
$data[23,3] = -0.5167593960592196
$data[23,4] = 0.1352448915017527
$data[23,5] = 0.9848229289054871
$data[23,6] = 0.9047340750694275
$data[23,7] = 0.3318836092948914
$data[23,8] = 0.6351562738418579

$data[24,0] = "model_2_1_0"
$data[24,1] = 0.112277992962789
$data[24,2] = 0.1197021215827284
$data[24,3] = -0.5793629665148985
$data[24,4] = 0.1522574885529031
$data[24,5] = 0.9824467301368713
$data[24,6] = 0.8689517974853516
$data[24,7] = 0.3455819487571716
$data[24,8] = 0.6226606965065002

$range = $ws.Range("A2:I26")
$range.Value = $data

Write-Output "Done"